$d = $word.ActiveDocument

# Locate the "Implementatie" heading and the last bullet of the feedback list
# ("Aantal mensen die hebben gestemd") by content, rather than by a fixed
# paragraph index, so the replacement range is found reliably.
$startPara = $null
$endPara = $null
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text.TrimEnd([char]13, [char]7)
    if ($t -eq "Implementatie") {
        $startPara = $p
    }
    if ($t -eq "Aantal mensen die hebben gestemd") {
        $endPara = $p
    }
}

$targetRange = $d.Range($startPara.Range.Start, $endPara.Range.End)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main" xmlns:wp="http://schemas.openxmlformats.org/drawingml/2006/wordprocessingDrawing" xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture" xmlns:r="http://schemas.openxmlformats.org/officeDocument/2006/relationships" xmlns:wp14="http://schemas.microsoft.com/office/word/2010/wordprocessingDrawing" xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main"><w:body><w:p><w:pPr><w:pStyle w:val="Heading1"/></w:pPr><w:proofErr w:type="gramStart"/><w:r><w:lastRenderedPageBreak/><w:t>Implementatie</w:t></w:r><w:proofErr w:type="gramEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="Heading3"/></w:pPr><w:r><w:t>Google Visualization</w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">In de </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:t>implementatie</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:t xml:space="preserve"> fase ben ik begonnen met het onderzoeken van Dimple.js en D3.js. Allereerst heb ik geprobeerd een simpele </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>scatter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> </w:t></w:r><w:r><w:t>plot</w:t></w:r><w:r><w:t xml:space="preserve"> te maken met Dimple. Echter liep ik steeds tegen obstakels aan en ben ik na een dag niets verder gekomen. Hierna heb ik geprobeerd om direct met D3 te werken maar zonder succes. </w:t></w:r></w:p><w:p><w:r><w:t>Om toch verder te kunnen met mijn project heb ik ervoor gekozen om van Dimple en D3 af te stappen en te zoeken naar een andere tool</w:t></w:r><w:r><w:t xml:space="preserve"> zodat ik me kan focussen op de visualisatie in plaats van het uitzoeken van deze technieken</w:t></w:r><w:r><w:t>. Al snel kwam ik tot Google Visualization</w:t></w:r><w:r><w:t xml:space="preserve">, hiermee heb ik erg snel een simpele </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:t>scatter</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:t xml:space="preserve"> plot kunnen maken. </w:t></w:r><w:r><w:t xml:space="preserve"> </w:t></w:r></w:p><w:p><w:r><w:t xml:space="preserve">Zo heb ik doorgewerkt naar het eerste prototype van mijn visualisatie die hieronder is te zien. </w:t></w:r><w:r><w:t>De volgende belangrijke functie die ik heb toegevoegd is een filter waarmee de gebruiker een selectie kan maken op basis van het jaar dat de film is uitgekomen. Met deze versie heb ik tijdens het klassikaal feedback moment de volgende feedback ontvangen:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Assen vast zetten</w:t></w:r><w:r><w:t xml:space="preserve"> om </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Titel toevoegen aan de website die de gebruiker uitdaagt om te onderzoeken.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t xml:space="preserve">Titel van film toevoegen </w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Naam van assen</w:t></w:r><w:r><w:t xml:space="preserve"> toevoegen</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:t>Mogelijk poster van film toevoege</w:t></w:r><w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/><w:r><w:t>n</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="ListParagraph"/><w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:noProof/><w:lang w:eastAsia="nl-NL"/></w:rPr><w:drawing><wp:anchor distT="0" distB="0" distL="114300" distR="114300" simplePos="0" relativeHeight="251664384" behindDoc="0" locked="0" layoutInCell="1" allowOverlap="1"><wp:simplePos x="0" y="0"/><wp:positionH relativeFrom="page"><wp:align>left</wp:align></wp:positionH><wp:positionV relativeFrom="paragraph"><wp:posOffset>433705</wp:posOffset></wp:positionV><wp:extent cx="7627620" cy="3825240"/><wp:effectExtent l="0" t="0" r="0" b="3810"/><wp:wrapThrough wrapText="bothSides"><wp:wrapPolygon edited="0"><wp:start x="0" y="0"/><wp:lineTo x="0" y="21514"/><wp:lineTo x="21524" y="21514"/><wp:lineTo x="21524" y="0"/><wp:lineTo x="0" y="0"/></wp:wrapPolygon></wp:wrapThrough><wp:docPr id="7" name="Picture 7" descr="https://gyazo.com/027bfa8af4b60bc835d4062b23a0df4e.png"/><wp:cNvGraphicFramePr><a:graphicFrameLocks xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main" noChangeAspect="1"/></wp:cNvGraphicFramePr><a:graphic xmlns:a="http://schemas.openxmlformats.org/drawingml/2006/main"><a:graphicData uri="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:pic xmlns:pic="http://schemas.openxmlformats.org/drawingml/2006/picture"><pic:nvPicPr><pic:cNvPr id="0" name="Picture 1" descr="https://gyazo.com/027bfa8af4b60bc835d4062b23a0df4e.png"/><pic:cNvPicPr><a:picLocks noChangeAspect="1" noChangeArrowheads="1"/></pic:cNvPicPr></pic:nvPicPr><pic:blipFill><a:blip r:embed="rId12"><a:extLst><a:ext uri="{28A0092B-C50C-407E-A947-70E740481C1C}"><a14:useLocalDpi xmlns:a14="http://schemas.microsoft.com/office/drawing/2010/main" val="0"/></a:ext></a:extLst></a:blip><a:srcRect/><a:stretch><a:fillRect/></a:stretch></pic:blipFill><pic:spPr bwMode="auto"><a:xfrm><a:off x="0" y="0"/><a:ext cx="7627620" cy="3825240"/></a:xfrm><a:prstGeom prst="rect"><a:avLst/></a:prstGeom><a:noFill/><a:ln><a:noFill/></a:ln></pic:spPr></pic:pic></a:graphicData></a:graphic><wp14:sizeRelH relativeFrom="margin"><wp14:pctWidth>0</wp14:pctWidth></wp14:sizeRelH><wp14:sizeRelV relativeFrom="margin"><wp14:pctHeight>0</wp14:pctHeight></wp14:sizeRelV></wp:anchor></w:drawing></w:r><w:r><w:t>Kijken of het a</w:t></w:r><w:r><w:t>antal mensen die hebben gestemd</w:t></w:r><w:r><w:t xml:space="preserve"> interessant is.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

[void]$targetRange.InsertXML($xml)

Write-Output "Paragraph count after edit: $($d.Paragraphs.Count)"
